$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same
$ws.Range("A1").Value = "Current Name"
$ws.Range("B1").Value = "New Name"

# Replace/extend data rows 2-8
$ws.Range("A2").Value = "gk-aks-Digital/firstgithubrepo"
$ws.Range("B2").Value = "gk-aks-Shared/firstgithubrepo"

$ws.Range("A3").Value = "gk-aks-Digital/secondgithubrepo"
$ws.Range("B3").Value = "gk-aks-Confidential/secondgithubrepo"

$ws.Range("A4").Value = "gk-aks-Digital/thirdgithubrepo"
$ws.Range("B4").Value = "gk-aks-Confidential/thirdgithubrepo"

$ws.Range("A5").Value = "gk-aks-Digital/fourthgithubrepo"
$ws.Range("B5").Value = "gk-aks-Shared/fourthgithubrepo"

$ws.Range("A6").Value = "gk-aks-Digital/fifthgithubrepo"
$ws.Range("B6").Value = "gk-aks-Confidential/fifthgithubrepo"

$ws.Range("A7").Value = "gk-aks-Digital/sixthgithubrepo"
$ws.Range("B7").Value = "gk-aks-Shared/sixthgithubrepo"

$ws.Range("A8").Value = "gk-aks-Digital/seventhgithubrepo"
$ws.Range("B8").Value = "gk-aks-Shared/seventhgithubrepo"

$ws.Range("B8").Select()
